# Insert a new data row above the current row 72, shifting all existing
# rows 72..96 down to 73..97 (dimension becomes A1:R97).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new record's data.
$ws.Cells.Item(72, 1).Value = 6
$ws.Cells.Item(72, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(72, 3).Value = "Metropolitana"
$ws.Cells.Item(72, 4).Value = 44460
$ws.Cells.Item(72, 5).Value = 13
$ws.Cells.Item(72, 6).Value = 100112001
$ws.Cells.Item(72, 7).Value = "Berenjena"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 310
$ws.Cells.Item(72, 11).Value = 7000
$ws.Cells.Item(72, 12).Value = 8000
$ws.Cells.Item(72, 13).Value = 7581
$ws.Cells.Item(72, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 152
$ws.Cells.Item(72, 17).Value = 50
$ws.Cells.Item(72, 18).Value = "Hortaliza"
